$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9588264299802761
$ws.Range("C2").Value = 0.02616596809291777

$ws.Range("B3").Value = 0.8855522682445759
$ws.Range("C3").Value = 0.03263349969368202

$ws.Range("B4").Value = 0.965483234714004
$ws.Range("C4").Value = 0.03851649843514515

$ws.Range("B5").Value = 0.8257889546351086
$ws.Range("C5").Value = 0.05407184155249235

$ws.Range("B6").Value = 0.9016642011834319
$ws.Range("C6").Value = 0.06465369704344144
